$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers): rename existing antibiotic columns to the "_1_" course
#     suffix, and insert two more Cefepime courses (course 2 and course 3)
#     ahead of Piperacillin/Tazobactam, which itself becomes course 1. ---

# Capture the original Piperacillin/Tazobactam data (currently in F:G) before
# those columns are overwritten with the new Cefepime course headers/values.
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2

$ws.Range("B1").Value2 = "Amoxicillin_1_FA"
$ws.Range("C1").Value2 = "Amoxicillin_1_LA"
$ws.Range("D1").Value2 = "Cefepime_1_FA"
$ws.Range("E1").Value2 = "Cefepime_1_LA"
$ws.Range("F1").Value2 = "Cefepime_2_FA"
$ws.Range("G1").Value2 = "Cefepime_2_LA"
$ws.Range("H1").Value2 = "Cefepime_3_FA"
$ws.Range("I1").Value2 = "Cefepime_3_LA"
$ws.Range("J1").Value2 = "Piperacillin/Tazobactam_1_FA"
$ws.Range("K1").Value2 = "Piperacillin/Tazobactam_1_LA"

# Match the existing header formatting (bold, centered, thin-bordered) on the
# newly added header cells H1:K1.
$ws.Range("H1:K1").Font.Bold = $true
$ws.Range("H1:K1").HorizontalAlignment = -4108
$ws.Range("H1:K1").VerticalAlignment = -4160
$ws.Range("H1:K1").Borders.LineStyle = 1

# --- Row 2 (data): relocate the original Piperacillin/Tazobactam values to
#     J2/K2, update the Cefepime course-1 last-admin time, and populate the
#     new course 2/3 timestamps. ---
$ws.Range("J2").Value2 = $oldF2
$ws.Range("K2").Value2 = $oldG2

$ws.Range("E2").Value2 = 44587.80813657407
$ws.Range("F2").Value2 = 44592.89236111111
$ws.Range("G2").Value2 = 44596.28888888889
$ws.Range("H2").Value2 = 44598.20555555556
$ws.Range("I2").Value2 = 44598.28888888889

# Apply the same date/time number format used by the other timestamp cells.
$ws.Range("F2:K2").NumberFormat = $ws.Range("D2").NumberFormat

# --- Row 3: extend the trailing blank cells through the new columns. ---
$ws.Range("H3").Formula = "=""""" 
$ws.Range("I3").Formula = "=""""" 
$ws.Range("J3").Formula = "=""""" 
$ws.Range("K3").Formula = "=""""" 
